$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header for column C: Obstacles -> obstacle_types
$ws.Range("C1").Value = "obstacle_types"

# Add new trial row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "Hallway"
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 0.5
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 64
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 0.99
$ws.Range("J9").Value = 4000
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 2
$ws.Range("M9").Value = "No"

# Move the selection cursor (matches author's final cursor position in the file)
$ws.Range("E17").Select()
